$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 08:07"

# Update Ucrania (row 28)
$ws.Range("B28").Value = 287231
$ws.Range("C28").Value = 5992
$ws.Range("D28").Value = 121919
$ws.Range("E28").Value = 159904
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 106
$ws.Range("H28").Value = 5408

# Update Uzbekistan (row 60)
$ws.Range("B60").Value = 62484
$ws.Range("C60").Value = 206
$ws.Range("D60").Value = 59429
$ws.Range("E60").Value = 2536
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 519

# Update Kirguistan (row 68)
$ws.Range("B68").Value = 51020
$ws.Range("C68").Value = 431
$ws.Range("D68").Value = 45288
$ws.Range("E68").Value = 4629
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 1103

# Update El Salvador (row 83)
$ws.Range("E83").Value = 3838
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 4
$ws.Range("H83").Value = 912
